$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 50-56: "Configurar perfil" / "Consulta de Usuarios" vote items.
# Column B uses center alignment (matches existing style index 2),
# column C uses left alignment (matches existing style index 3),
# column D ("votar") keeps the default style.

$rows = @(
    @{ Row = 50; B = "Consulta de Usuarios"; C = "Eu, ADM quero alterar status do usuário." },
    @{ Row = 51; B = "Consulta de Usuarios"; C = "Eu, ADM quero acessar a tela de editar usuário" },
    @{ Row = 52; B = "Consulta de Usuarios"; C = "Eu, ADM quero redefinir senha do usuário" },
    @{ Row = 53; B = "Consulta de Usuarios"; C = "Eu, ADM quero salvar alterações no usuário" },
    @{ Row = 54; B = "Configurar perfil";    C = "Eu, ADM quero acessar tela de Permissões de perfis" },
    @{ Row = 55; B = "Configurar perfil";    C = "Eu, ADM quero fazer alterações nas permissões dos perfis" },
    @{ Row = 56; B = "Configurar perfil";    C = "Eu, ADM quero salvar alterações nas permissões dos perfis" }
)

foreach ($r in $rows) {
    $cCell = $ws.Cells.Item($r.Row, 3)
    $cCell.Value = $r.C
    $cCell.HorizontalAlignment = -4131
}

foreach ($r in $rows) {
    $bCell = $ws.Cells.Item($r.Row, 2)
    $bCell.Value = $r.B
    $bCell.HorizontalAlignment = -4108
}

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value = "votar"
}

$ws.Range("D56").Select()
